# Insert a new parameter row for "chemical_recycling_pyrolysis" right after
# the existing "chemical_recycling_gasification" row (row 9), pushing every
# row below it (fossil_routes ... fossil_lock_in) down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("10:10").Insert()

$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
